# Fix Training Data Issue: the "Date" column (BF) was off by one day due to
# how NBA stats were shown. Correct every row's BF value from the old
# "6-10-2013-14" label to the proper ISO date string "2014-06-10".
# The leading apostrophe forces Excel to keep the value as literal text
# instead of re-interpreting the date-shaped string as a serial date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BF2:BF31").Value = "'2014-06-10"
